$wb = $excel.ActiveWorkbook

$status = "Handback transform failed"

# "Status" column (row 3 = the 9f199e39-... file) moves from "Ready for
# handoff" to "Handback transform failed" everywhere it is shown: the
# Overview roll-up (both the zh-cn and de-de status columns) and each
# language sheet's own Status column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status

# New "Error Detail" entries explaining the handback mismatch, one per
# language sheet, for the same row.
$wsZhCn.Range("L3").Value = "Handback file name: ayqjl4ui.5ot is different with handoff file name: 9f199e39-1be8-44f6-8a3f-19d3f101795a.00d6cd9f64a91688e32da06d3fbbc52573328c1f.zh-cn."
$wsDeDe.Range("L3").Value = "Handback file name: ayqjl4ui.5ot is different with handoff file name: 9f199e39-1be8-44f6-8a3f-19d3f101795a.00d6cd9f64a91688e32da06d3fbbc52573328c1f.de-de."
